$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New column I (rows 1-11): mirror the existing column J literal values.
#    (No explicit style -> default style 0, unlike J which already carries
#    style 2/3.)
# ---------------------------------------------------------------------------
$ws.Range("I1").Value = 62021
$ws.Range("I2").Value = 126951
$ws.Range("I3").Value = 12471
$ws.Range("I4").Value = 662
$ws.Range("I5").Value = 124991
$ws.Range("I6").Value = 78589
$ws.Range("I7").Value = 5115
$ws.Range("I8").Value = 53052
$ws.Range("I9").Value = 105928
$ws.Range("I10").Value = 948
$ws.Range("I11").Value = 1719

# ---------------------------------------------------------------------------
# 2. New rows 15-25: vote-transfer table for the 2019 NI European
#    Parliament election count, plus a trailing styled-but-empty row 26.
#    Column A uses the "Surname, First" shared strings already present in
#    the workbook (indices 12-22); B/D/F/H/J are literal running totals,
#    C/E/G/I are the per-stage deltas (formulas where the source workbook
#    used one).
# ---------------------------------------------------------------------------

# Row 15 - Jim Allister
$ws.Range("A15").Value = "Allister, Jim"
$ws.Range("B15").Value = 62021
$ws.Range("C15").Formula = "=D15-B15"
$ws.Range("D15").Value = 63872
$ws.Range("E15").Formula = "=F15-D15"
$ws.Range("F15").Value = 79540
$ws.Range("G15").Formula = "=H15-F15"
$ws.Range("H15").Value = 89854
$ws.Range("I15").Formula = "=J15-H15"
$ws.Range("J15").Value = 90079

# Row 16 - Martina Anderson
$ws.Range("A16").Value = "Anderson, Martina"
$ws.Range("B16").Value = 126951
$ws.Range("C16").Formula = "=D16-B16"
$ws.Range("D16").Value = 128117
$ws.Range("E16").Formula = "=F16-D16"
$ws.Range("F16").Value = 127190
$ws.Range("G16").Formula = "=H16-F16"
$ws.Range("H16").Value = 128200.5
$ws.Range("I16").Formula = "=J16-H16"
$ws.Range("J16").Value = 152436.5

# Row 17 - Clare Bailey
$ws.Range("A17").Value = "Bailey, Clare"
$ws.Range("B17").Value = 12471
$ws.Range("C17").Formula = "=-B17"
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0

# Row 18 - Amandeep Bhogal
$ws.Range("A18").Value = "Bhogal, Amandeep"
$ws.Range("B18").Value = 662
$ws.Range("C18").Formula = "=-B18"
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0

# Row 19 - Diane Dodds
$ws.Range("A19").Value = "Dodds, Diane"
$ws.Range("B19").Value = 124991
$ws.Range("C19").Formula = "=D19-B19"
$ws.Range("D19").Value = 127291
$ws.Range("E19").Formula = "=F19-D19"
$ws.Range("F19").Value = 155422
$ws.Range("G19").Formula = "=143112-F19"
$ws.Range("H19").Value = 143112
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 143112

# Row 20 - Colum Eastwood (no C20 cell in the source)
$ws.Range("A20").Value = "Eastwood, Colum"
$ws.Range("B20").Value = 78589
$ws.Range("D20").Value = 80949
$ws.Range("E20").Formula = "=F20-D20"
$ws.Range("F20").Value = 82101
$ws.Range("G20").Formula = "=H20-F20"
$ws.Range("H20").Value = 82413.5
$ws.Range("I20").Formula = "=-H20"
$ws.Range("J20").Value = 0

# Row 21 - Robert Hill
$ws.Range("A21").Value = "Hill, Robert"
$ws.Range("B21").Value = 5115
$ws.Range("C21").Formula = "=-B21"
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0

# Row 22 - Danny Kennedy (no C22 cell in the source)
$ws.Range("A22").Value = "Kennedy, Danny"
$ws.Range("B22").Value = 53052
$ws.Range("D22").Value = 53052
$ws.Range("E22").Formula = "=-D22"
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0

# Row 23 - Naomi Long
$ws.Range("A23").Value = "Long, Naomi"
$ws.Range("B23").Value = 105928
$ws.Range("C23").Formula = "=D23-B23"
$ws.Range("D23").Value = 115327
$ws.Range("E23").Formula = "=F23-D23"
$ws.Range("F23").Value = 122263
$ws.Range("G23").Formula = "=H23-F23"
$ws.Range("H23").Value = 123917
$ws.Range("I23").Formula = "=J23-H23"
$ws.Range("J23").Value = 170370

# Row 24 - Neil McCann
$ws.Range("A24").Value = "McCann, Neil"
$ws.Range("B24").Value = 948
$ws.Range("C24").Formula = "=-B24"
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0

# Row 25 - Jane Morrice
$ws.Range("A25").Value = "Morrice, Jane"
$ws.Range("B25").Value = 1719
$ws.Range("C25").Formula = "=-B25"
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0

# Row 26 - trailing styled-but-empty cell so the dimension stretches to J26
$ws.Range("J26").Value = ""

# ---------------------------------------------------------------------------
# 3. Styling: H15:H25 and J15:J26 (and I20) carry a new 2-decimal number
#    format (numFmtId 2, "0.00") that doesn't exist yet in the workbook -
#    setting NumberFormat here appends a brand-new cellXfs entry exactly
#    like the diff (cellXfs count 5 -> 6).
# ---------------------------------------------------------------------------
$ws.Range("H15:H25").NumberFormat = "0.00"
$ws.Range("J15:J26").NumberFormat = "0.00"
$ws.Range("I20").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 4. Column widths for the new I and J header cells, dimension/view refresh.
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 8.6
$ws.Columns.Item(10).ColumnWidth = 8.6

# ---------------------------------------------------------------------------
# 5. Selection / view state, matching the post-edit workbook (scrolled down
#    to the newly-added rows, C15:J25 selected).
# ---------------------------------------------------------------------------
$ws.Range("C15:J25").Select()
$excel.ActiveWindow.ScrollRow = 10
